$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
# Row 34
$ws.Range("H34").Value = 11400
$ws.Range("I34").Value = 11400
$ws.Range("K34").Value = 11400
$ws.Range("M34").Value = -11197

# Row 36
$ws.Range("H36").Value = 11400
$ws.Range("I36").Value = 11400
$ws.Range("K36").Value = 11400
$ws.Range("M36").Value = -10685

# Row 76
$ws.Range("H76").Value = 3906404
$ws.Range("I76").Value = 4686985
$ws.Range("K76").Value = 4686985
$ws.Range("M76").Value = -4686670

# Row 79
$ws.Range("H79").Value = 3906404
$ws.Range("I79").Value = 4686985
$ws.Range("K79").Value = 4686985
$ws.Range("M79").Value = -4685893

# Row 88
$ws.Range("H88").Value = 55557920
$ws.Range("I88").Value = 166667900
$ws.Range("K88").Value = 166667900
$ws.Range("M88").Value = -166667494

# Row 91
$ws.Range("H91").Value = 55557920
$ws.Range("I91").Value = 166667900
$ws.Range("K91").Value = 166667900
$ws.Range("M91").Value = -166666496

# Row 100
$ws.Range("H100").Value = 1513.75
$ws.Range("I100").Value = 925
$ws.Range("K100").Value = 925
$ws.Range("M100").Value = -384

# Row 127
$ws.Range("H127").Value = 2291.4119
$ws.Range("I127").Value = 1996.5
$ws.Range("K127").Value = 5989.5
$ws.Range("M127").Value = -1029.5

# Row 129
$ws.Range("H129").Value = 881.61816
$ws.Range("I129").Value = 491.77777
$ws.Range("K129").Value = 1475.33331
$ws.Range("M129").Value = 3524.66669

# Row 132
$ws.Range("H132").Value = 4976121.5
$ws.Range("I132").Value = 5556478
$ws.Range("K132").Value = 16669434
$ws.Range("M132").Value = -16666904

# Row 137
$ws.Range("H137").Value = 1053922
$ws.Range("I137").Value = 1057.2727
$ws.Range("J137").Value = 2501611
$ws.Range("K137").Value = 3171.8181
$ws.Range("L137").Value = 7504833
$ws.Range("M137").Value = -621.8181
$ws.Range("N137").Value = -7509933

# Row 138
$ws.Range("H138").Value = 1609.1464
$ws.Range("I138").Value = 1299.8474
$ws.Range("J138").Value = 2402.5652
$ws.Range("K138").Value = 3899.5422
$ws.Range("L138").Value = 7207.6956
$ws.Range("M138").Value = 1240.4578
$ws.Range("N138").Value = -17487.6956

# Row 139
$ws.Range("H139").Value = 51828.4
$ws.Range("J139").Value = 51828.4
$ws.Range("L139").Value = 51828.4
$ws.Range("N139").Value = -62108.4

# Row 141
$ws.Range("H141").Value = 849873.5
$ws.Range("I141").Value = 966276.8
$ws.Range("J141").Value = 5949.75
$ws.Range("K141").Value = 2898830.4
$ws.Range("L141").Value = 17849.25
$ws.Range("M141").Value = -2893650.4
$ws.Range("N141").Value = -28209.25


# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 4123.9854
$ws.Range("I32").Value = 3339.1147
$ws.Range("J32").Value = 10963.571
$ws.Range("K32").Value = 3339.1147
$ws.Range("L32").Value = 10963.571
$ws.Range("M32").Value = -3052.1147
$ws.Range("N32").Value = -11537.571

# Row 63
$ws.Range("H63").Value = 8143.5713
$ws.Range("I63").Value = 8101
$ws.Range("J63").Value = 8250
$ws.Range("K63").Value = 8101
$ws.Range("L63").Value = 8250
$ws.Range("M63").Value = -7415
$ws.Range("N63").Value = -9622

# Row 66
$ws.Range("H66").Value = 8143.5713
$ws.Range("I66").Value = 8101
$ws.Range("J66").Value = 8250
$ws.Range("K66").Value = 40505
$ws.Range("L66").Value = 41250
$ws.Range("M66").Value = -37073
$ws.Range("N66").Value = -48114

# Row 122
$ws.Range("H122").Value = 1559.9269
$ws.Range("I122").Value = 1214.0938
$ws.Range("J122").Value = 2789.5557
$ws.Range("K122").Value = 3642.2814
$ws.Range("L122").Value = 8368.667099999999
$ws.Range("M122").Value = -1192.2814
$ws.Range("N122").Value = -13268.6671

# Row 132
$ws.Range("H132").Value = 1326.1372
$ws.Range("I132").Value = 1043.3846
$ws.Range("J132").Value = 2245.0833
$ws.Range("K132").Value = 3130.1538
$ws.Range("L132").Value = 6735.249899999999
$ws.Range("M132").Value = -600.1538
$ws.Range("N132").Value = -11795.2499


# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
# Row 86
$ws.Range("H86").Value = 416885.34
$ws.Range("I86").Value = 718116.5
$ws.Range("K86").Value = 718116.5
$ws.Range("M86").Value = -716993.5

# Row 89
$ws.Range("H89").Value = 416885.34
$ws.Range("I89").Value = 718116.5
$ws.Range("K89").Value = 3590582.5
$ws.Range("M89").Value = -3584966.5

# Row 134
$ws.Range("H134").Value = 1247.4
$ws.Range("I134").Value = 1274.2941
$ws.Range("J134").Value = 333
$ws.Range("K134").Value = 3822.8823
$ws.Range("L134").Value = 999
$ws.Range("M134").Value = -1287.8823
$ws.Range("N134").Value = -6069


# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
# Row 99
$ws.Range("H99").Value = 1132.75
$ws.Range("I99").Value = 1132.75
$ws.Range("K99").Value = 1132.75
$ws.Range("M99").Value = 365.25

# Row 107
$ws.Range("H107").Value = 1243.8572
$ws.Range("I107").Value = 1243.8572
$ws.Range("K107").Value = 1243.8572
$ws.Range("M107").Value = 676.1428000000001

# Row 122
$ws.Range("H122").Value = 3886.2354
$ws.Range("I122").Value = 2711.2
$ws.Range("J122").Value = 5564.857
$ws.Range("K122").Value = 8133.599999999999
$ws.Range("L122").Value = 16694.571
$ws.Range("M122").Value = -5683.599999999999
$ws.Range("N122").Value = -21594.571

# Row 126
$ws.Range("H126").Value = 1132.75
$ws.Range("I126").Value = 1132.75
$ws.Range("K126").Value = 3398.25
$ws.Range("M126").Value = -928.25

# Row 132
$ws.Range("H132").Value = 1710.6364
$ws.Range("I132").Value = 1139.1666
$ws.Range("J132").Value = 3234.5557
$ws.Range("K132").Value = 3417.4998
$ws.Range("L132").Value = 9703.667099999999
$ws.Range("M132").Value = -887.4998000000001
$ws.Range("N132").Value = -14763.6671

# Row 141
$ws.Range("H141").Value = 66989
$ws.Range("J141").Value = 66989
$ws.Range("L141").Value = 66989
$ws.Range("N141").Value = -77349


# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
# Row 32
$ws.Range("H32").Value = 966.6667
$ws.Range("J32").Value = 966.6667
$ws.Range("L32").Value = 2900.0001
$ws.Range("N32").Value = -3466.0001

# Row 81
$ws.Range("H81").Value = 2296
$ws.Range("I81").Value = 1980
$ws.Range("J81").Value = 2359.2
$ws.Range("K81").Value = 5940
$ws.Range("L81").Value = 7077.599999999999
$ws.Range("N81").Value = -9323.599999999999
$ws.Range("M81").Value = -4817

# Row 84
$ws.Range("H84").Value = 2296
$ws.Range("I84").Value = 1980
$ws.Range("J84").Value = 2359.2
$ws.Range("K84").Value = 17820
$ws.Range("L84").Value = 21232.8
$ws.Range("N84").Value = -32464.8
$ws.Range("M84").Value = -12204

# Row 131
$ws.Range("H131").Value = 5271373
$ws.Range("J131").Value = 9331.253000000001
$ws.Range("L131").Value = 27993.759
$ws.Range("N131").Value = -38073.75900000001

# Row 137
$ws.Range("H137").Value = 3053.2144
$ws.Range("J137").Value = 4254.353
$ws.Range("L137").Value = 12763.059
$ws.Range("N137").Value = -22963.059


# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
# Row 102
$ws.Range("H102").Value = 2785.6667
$ws.Range("I102").Value = 2811.6365
$ws.Range("K102").Value = 2811.6365
$ws.Range("M102").Value = -1189.6365

# Row 122
$ws.Range("H122").Value = 1431.1666
$ws.Range("I122").Value = 1147.75
$ws.Range("K122").Value = 3443.25
$ws.Range("M122").Value = -993.25

# Row 126
$ws.Range("H126").Value = 2418207.2
$ws.Range("I126").Value = 9262585
$ws.Range("J126").Value = 2544.5881
$ws.Range("K126").Value = 27787755
$ws.Range("L126").Value = 7633.7643
$ws.Range("M126").Value = -27785285
$ws.Range("N126").Value = -12573.7643

# Row 132
$ws.Range("H132").Value = 1481554.5
$ws.Range("I132").Value = 1833201.9
$ws.Range("J132").Value = 4635.6
$ws.Range("K132").Value = 5499605.699999999
$ws.Range("L132").Value = 13906.8
$ws.Range("M132").Value = -5497075.699999999
$ws.Range("N132").Value = -18966.8


# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
# Row 22
$ws.Range("H22").Value = 2317
$ws.Range("I22").Value = 10000
$ws.Range("K22").Value = 10000
$ws.Range("M22").Value = -9705

# Row 27
$ws.Range("H27").Value = 2317
$ws.Range("I27").Value = 10000
$ws.Range("K27").Value = 10000
$ws.Range("M27").Value = -9893

# Row 132
$ws.Range("H132").Value = 1128.7971
$ws.Range("I132").Value = 914.6531
$ws.Range("K132").Value = 2743.9593
$ws.Range("M132").Value = -213.9593

# Row 136
$ws.Range("H136").Value = 1649.9824
$ws.Range("I136").Value = 1020.70215
$ws.Range("J136").Value = 4607.6
$ws.Range("K136").Value = 3062.10645
$ws.Range("L136").Value = 13822.8
$ws.Range("M136").Value = -512.1064499999998
$ws.Range("N136").Value = -18922.8


# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
# Row 14
$ws.Range("H14").Value = 7045.8335
$ws.Range("I14").Value = 6989.3335
$ws.Range("J14").Value = 7074.0835
$ws.Range("K14").Value = 6989.3335
$ws.Range("L14").Value = 7074.0835
$ws.Range("M14").Value = -6821.3335
$ws.Range("N14").Value = -7410.0835

# Row 81
$ws.Range("H81").Value = 716.6667
$ws.Range("I81").Value = 400.5
$ws.Range("J81").Value = 874.75
$ws.Range("K81").Value = 801
$ws.Range("L81").Value = 1749.5
$ws.Range("M81").Value = 260
$ws.Range("N81").Value = -3871.5

# Row 84
$ws.Range("H84").Value = 716.6667
$ws.Range("I84").Value = 400.5
$ws.Range("J84").Value = 874.75
$ws.Range("K84").Value = 4005
$ws.Range("L84").Value = 8747.5
$ws.Range("M84").Value = 1299
$ws.Range("N84").Value = -19355.5

# Row 113
$ws.Range("H113").Value = 644.4
$ws.Range("I113").Value = 305.58334
$ws.Range("K113").Value = 916.7500200000001
$ws.Range("M113").Value = 1253.24998

# Row 122
$ws.Range("H122").Value = 31490.576
$ws.Range("I122").Value = 32686.2
$ws.Range("K122").Value = 98058.60000000001
$ws.Range("M122").Value = -95608.60000000001

# Row 126
$ws.Range("H126").Value = 6756.4585
$ws.Range("I126").Value = 10738.091
$ws.Range("J126").Value = 3387.3845
$ws.Range("K126").Value = 32214.273
$ws.Range("L126").Value = 10162.1535
$ws.Range("M126").Value = -29744.273
$ws.Range("N126").Value = -15102.1535

# Row 132
$ws.Range("H132").Value = 1342.2322
$ws.Range("I132").Value = 1017.2292
$ws.Range("J132").Value = 3292.25
$ws.Range("K132").Value = 3051.6876
$ws.Range("L132").Value = 9876.75
$ws.Range("M132").Value = -521.6876000000002
$ws.Range("N132").Value = -14936.75

# Row 136
$ws.Range("H136").Value = 13229046
$ws.Range("I136").Value = 14621209
$ws.Range("K136").Value = 43863627
$ws.Range("M136").Value = -43861077

